# Auto-generated edit script applying the scheduled-runner market data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1192.6666
$ws.Range("I41").Value = 92.333336
$ws.Range("K41").Value = 92.333336
$ws.Range("M41").Value = 347.666664
$ws.Range("H62").Value = 70373200
$ws.Range("I62").Value = 90478776
$ws.Range("J62").Value = 3694.75
$ws.Range("K62").Value = 90478776
$ws.Range("L62").Value = 3694.75
$ws.Range("M62").Value = -90478152
$ws.Range("N62").Value = -4942.75
$ws.Range("H65").Value = 70373200
$ws.Range("I65").Value = 90478776
$ws.Range("J65").Value = 3694.75
$ws.Range("K65").Value = 452393880
$ws.Range("L65").Value = 18473.75
$ws.Range("M65").Value = -452390760
$ws.Range("N65").Value = -24713.75
$ws.Range("H113").Value = 79416.234
$ws.Range("I113").Value = 2928.4285
$ws.Range("J113").Value = 168652
$ws.Range("K113").Value = 2928.4285
$ws.Range("L113").Value = 168652
$ws.Range("M113").Value = 325.5715
$ws.Range("N113").Value = -175160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 382.5
$ws.Range("I5").Value = 516.6667
$ws.Range("J5").Value = 248.33333
$ws.Range("K5").Value = 516.6667
$ws.Range("L5").Value = 248.33333
$ws.Range("M5").Value = -404.6667
$ws.Range("N5").Value = -472.33333
$ws.Range("H32").Value = 4881.029
$ws.Range("I32").Value = 5142.606
$ws.Range("K32").Value = 5142.606
$ws.Range("M32").Value = -4855.606
$ws.Range("H110").Value = 18751618
$ws.Range("I110").Value = 21429938
$ws.Range("K110").Value = 21429938
$ws.Range("M110").Value = -21427893

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 382.5
$ws.Range("I4").Value = 516.6667
$ws.Range("J4").Value = 248.33333
$ws.Range("K4").Value = 516.6667
$ws.Range("L4").Value = 248.33333
$ws.Range("M4").Value = -401.6667
$ws.Range("N4").Value = -478.33333
$ws.Range("H20").Value = 3049.7856
$ws.Range("I20").Value = 2643.8
$ws.Range("J20").Value = 4064.75
$ws.Range("K20").Value = 2643.8
$ws.Range("L20").Value = 4064.75
$ws.Range("M20").Value = -2396.8
$ws.Range("N20").Value = -4558.75
$ws.Range("H105").Value = 52645336
$ws.Range("I105").Value = 62516250
$ws.Range("K105").Value = 62516250
$ws.Range("M105").Value = -62514503

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 736.3333
$ws.Range("J8").Value = 1002.5
$ws.Range("L8").Value = 1002.5
$ws.Range("N8").Value = -1282.5
$ws.Range("H15").Value = 14
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 14
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H16").Value = 4347.95
$ws.Range("I16").Value = 3697.077
$ws.Range("K16").Value = 3697.077
$ws.Range("M16").Value = -3410.077
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9826
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H26").Value = 89.75
$ws.Range("J26").Value = 116.333336
$ws.Range("L26").Value = 116.333336
$ws.Range("N26").Value = -690.333336
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null
$ws.Range("H29").Value = 6281.8945
$ws.Range("J29").Value = 6463.1665
$ws.Range("L29").Value = 6463.1665
$ws.Range("N29").Value = -7049.1665
$ws.Range("H99").Value = 3158.7273
$ws.Range("I99").Value = 3158.7273
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3158.7273
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
$ws.Range("H107").Value = 1043.0588
$ws.Range("I107").Value = 818.3
$ws.Range("K107").Value = 818.3
$ws.Range("M107").Value = 1101.7
$ws.Range("H113").Value = 4347.95
$ws.Range("I113").Value = 3697.077
$ws.Range("K113").Value = 3697.077
$ws.Range("M113").Value = -1527.077
$ws.Range("H122").Value = 59787.47
$ws.Range("I122").Value = 63467.938
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 190403.814
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -187953.814
$ws.Range("N122").Value = -7600
$ws.Range("H126").Value = 3158.7273
$ws.Range("I126").Value = 3158.7273
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9476.1819
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 123698.4
$ws.Range("I132").Value = 2560.5
$ws.Range("K132").Value = 7681.5
$ws.Range("M132").Value = -5151.5
$ws.Range("H134").Value = 7726.7666
$ws.Range("I134").Value = 7638.615
$ws.Range("K134").Value = 22915.845
$ws.Range("M134").Value = -20380.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1599
$ws.Range("J45").Value = 1599
$ws.Range("L45").Value = 4797
$ws.Range("N45").Value = -5861
$ws.Range("H88").Value = 3999.1428
$ws.Range("I88").Value = 3999
$ws.Range("K88").Value = 11997
$ws.Range("M88").Value = -11569
$ws.Range("H91").Value = 3999.1428
$ws.Range("I91").Value = 3999
$ws.Range("K91").Value = 11997
$ws.Range("M91").Value = -10515
$ws.Range("H98").Value = 300
$ws.Range("I98").Value = 300
$ws.Range("K98").Value = 900
$ws.Range("M98").Value = 598
$ws.Range("H113").Value = 3509.9333
$ws.Range("I113").Value = 4079.9
$ws.Range("J113").Value = 2370
$ws.Range("K113").Value = 12239.7
$ws.Range("L113").Value = 7110
$ws.Range("M113").Value = -10069.7
$ws.Range("N113").Value = -11450

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 27666.666
$ws.Range("J69").Value = 27666.666
$ws.Range("L69").Value = 27666.666
$ws.Range("N69").Value = -29164.666
$ws.Range("H70").Value = 10172.32
$ws.Range("I70").Value = 9240.111000000001
$ws.Range("J70").Value = 12569.429
$ws.Range("K70").Value = 9240.111000000001
$ws.Range("L70").Value = 12569.429
$ws.Range("M70").Value = -8970.111000000001
$ws.Range("N70").Value = -13109.429
$ws.Range("H72").Value = 27666.666
$ws.Range("J72").Value = 27666.666
$ws.Range("L72").Value = 82999.99800000001
$ws.Range("N72").Value = -90487.99800000001
$ws.Range("H73").Value = 10172.32
$ws.Range("I73").Value = 9240.111000000001
$ws.Range("J73").Value = 12569.429
$ws.Range("K73").Value = 9240.111000000001
$ws.Range("L73").Value = 12569.429
$ws.Range("M73").Value = -8304.111000000001
$ws.Range("N73").Value = -14441.429
$ws.Range("H102").Value = 1400.7368
$ws.Range("I102").Value = 1087.2307
$ws.Range("J102").Value = 2080
$ws.Range("K102").Value = 1087.2307
$ws.Range("L102").Value = 2080
$ws.Range("M102").Value = 534.7692999999999
$ws.Range("N102").Value = -5324
$ws.Range("H132").Value = 4927.5
$ws.Range("I132").Value = 2745.7693
$ws.Range("J132").Value = 10600
$ws.Range("K132").Value = 8237.3079
$ws.Range("L132").Value = 31800
$ws.Range("M132").Value = -5707.3079
$ws.Range("N132").Value = -36860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 35723216
$ws.Range("J40").Value = 3675
$ws.Range("L40").Value = 3675
$ws.Range("N40").Value = -3947
$ws.Range("H68").Value = 2116.7693
$ws.Range("I68").Value = 2015.8334
$ws.Range("J68").Value = 2203.2856
$ws.Range("K68").Value = 2015.8334
$ws.Range("L68").Value = 2203.2856
$ws.Range("M68").Value = -1266.8334
$ws.Range("N68").Value = -3701.2856
$ws.Range("H71").Value = 2116.7693
$ws.Range("I71").Value = 2015.8334
$ws.Range("J71").Value = 2203.2856
$ws.Range("K71").Value = 10079.167
$ws.Range("L71").Value = 11016.428
$ws.Range("M71").Value = -6335.166999999999
$ws.Range("N71").Value = -18504.428
$ws.Range("H93").Value = 3201.9524
$ws.Range("I93").Value = 1983.0834
$ws.Range("J93").Value = 4827.1113
$ws.Range("K93").Value = 1983.0834
$ws.Range("L93").Value = 4827.1113
$ws.Range("M93").Value = -735.0834
$ws.Range("N93").Value = -7323.1113
$ws.Range("H124").Value = 64214.5
$ws.Range("J124").Value = 64214.5
$ws.Range("L124").Value = 64214.5
$ws.Range("N124").Value = -74034.5
$ws.Range("H129").Value = 81788
$ws.Range("J129").Value = 52750
$ws.Range("L129").Value = 52750
$ws.Range("N129").Value = -62750
$ws.Range("H136").Value = 4717.185
$ws.Range("I136").Value = 3465.5557
$ws.Range("K136").Value = 10396.6671
$ws.Range("M136").Value = -7846.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1970.875
$ws.Range("I107").Value = 2119.5
$ws.Range("J107").Value = 1525
$ws.Range("K107").Value = 6358.5
$ws.Range("L107").Value = 4575
$ws.Range("M107").Value = -4438.5
$ws.Range("N107").Value = -8415
